$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 904.79034
$ws.Range("J129").Value = 910.7692
$ws.Range("L129").Value = 2732.3076
$ws.Range("N129").Value = -12732.3076

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 23725.834
$ws.Range("J24").Value = 23725.834
$ws.Range("L24").Value = 23725.834
$ws.Range("N24").Value = -24473.834

$ws.Range("H61").Value = 1312.15
$ws.Range("I61").Value = 1359.619
$ws.Range("J61").Value = 1259.6842
$ws.Range("K61").Value = 1359.619
$ws.Range("L61").Value = 1259.6842
$ws.Range("M61").Value = -1147.619
$ws.Range("N61").Value = -1683.6842

$ws.Range("H74").Value = 9616445
$ws.Range("I74").Value = 11629058
$ws.Range("J74").Value = 629.2222
$ws.Range("K74").Value = 11629058
$ws.Range("L74").Value = 629.2222
$ws.Range("M74").Value = -11628184
$ws.Range("N74").Value = -2377.2222

$ws.Range("H77").Value = 9616445
$ws.Range("I77").Value = 11629058
$ws.Range("J77").Value = 629.2222
$ws.Range("K77").Value = 58145290
$ws.Range("L77").Value = 3146.111
$ws.Range("M77").Value = -58140922
$ws.Range("N77").Value = -11882.111

$ws.Range("H100").Value = 23725.834
$ws.Range("J100").Value = 23725.834
$ws.Range("L100").Value = 23725.834
$ws.Range("N100").Value = -25889.834

$ws.Range("H110").Value = 1731.1111
$ws.Range("I110").Value = 657.1429000000001
$ws.Range("J110").Value = 5490
$ws.Range("K110").Value = 657.1429000000001
$ws.Range("L110").Value = 5490
$ws.Range("M110").Value = 1387.8571
$ws.Range("N110").Value = -9580

$ws.Range("H136").Value = 1312.15
$ws.Range("I136").Value = 1359.619
$ws.Range("J136").Value = 1259.6842
$ws.Range("K136").Value = 4078.857
$ws.Range("L136").Value = 3779.0526
$ws.Range("M136").Value = -1528.857
$ws.Range("N136").Value = -8879.052599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2587102
$ws.Range("I86").Value = 3500
$ws.Range("J86").Value = 4653984
$ws.Range("K86").Value = 3500
$ws.Range("L86").Value = 4653984
$ws.Range("M86").Value = -2377
$ws.Range("N86").Value = -4656230

$ws.Range("H89").Value = 2587102
$ws.Range("I89").Value = 3500
$ws.Range("J89").Value = 4653984
$ws.Range("K89").Value = 17500
$ws.Range("L89").Value = 23269920
$ws.Range("M89").Value = -11884
$ws.Range("N89").Value = -23281152

$ws.Range("H94").Value = 846.0909
$ws.Range("I94").Value = 811.8889
$ws.Range("J94").Value = 1000
$ws.Range("K94").Value = 811.8889
$ws.Range("L94").Value = 1000
$ws.Range("M94").Value = -360.8889
$ws.Range("N94").Value = -1902

$ws.Range("H134").Value = 2852909.8
$ws.Range("I134").Value = 836.0833
$ws.Range("J134").Value = 7416227.5
$ws.Range("K134").Value = 2508.2499
$ws.Range("L134").Value = 22248682.5
$ws.Range("M134").Value = 26.7501000000002
$ws.Range("N134").Value = -22253752.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 15625830
$ws.Range("I122").Value = 25000798
$ws.Range("J122").Value = 884.1667
$ws.Range("K122").Value = 75002394
$ws.Range("L122").Value = 2652.5001
$ws.Range("M122").Value = -74999944
$ws.Range("N122").Value = -7552.5001

$ws.Range("H134").Value = 883.8444
$ws.Range("I134").Value = 817.3939
$ws.Range("J134").Value = 1066.5834
$ws.Range("K134").Value = 2452.1817
$ws.Range("L134").Value = 3199.7502
$ws.Range("M134").Value = 82.81829999999991
$ws.Range("N134").Value = -8269.7502

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 35722440
$ws.Range("I122").Value = 71428984
$ws.Range("J122").Value = 15900.571
$ws.Range("K122").Value = 642860856
$ws.Range("L122").Value = 143105.139
$ws.Range("M122").Value = -642858406
$ws.Range("N122").Value = -148005.139

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1513.421
$ws.Range("I102").Value = 1499.5454
$ws.Range("J102").Value = 1532.5
$ws.Range("K102").Value = 1499.5454
$ws.Range("L102").Value = 1532.5
$ws.Range("M102").Value = 122.4546
$ws.Range("N102").Value = -4776.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1900.5
$ws.Range("I7").Value = 1801
$ws.Range("J7").Value = 2000
$ws.Range("K7").Value = 1801
$ws.Range("L7").Value = 2000
$ws.Range("M7").Value = -1689
$ws.Range("N7").Value = -2224

$ws.Range("H40").Value = 16667989
$ws.Range("I40").Value = 1199.909
$ws.Range("J40").Value = 62501660
$ws.Range("K40").Value = 1199.909
$ws.Range("L40").Value = 62501660
$ws.Range("M40").Value = -1063.909
$ws.Range("N40").Value = -62501932

$ws.Range("H68").Value = 1226.2693
$ws.Range("I68").Value = 1135.9474
$ws.Range("J68").Value = 1471.4286
$ws.Range("K68").Value = 1135.9474
$ws.Range("L68").Value = 1471.4286
$ws.Range("M68").Value = -386.9474
$ws.Range("N68").Value = -2969.4286

$ws.Range("H71").Value = 1226.2693
$ws.Range("I71").Value = 1135.9474
$ws.Range("J71").Value = 1471.4286
$ws.Range("K71").Value = 5679.737
$ws.Range("L71").Value = 7357.143
$ws.Range("M71").Value = -1935.737
$ws.Range("N71").Value = -14845.143

$ws.Range("H122").Value = 5488.5557
$ws.Range("I122").Value = 7172.905
$ws.Range("J122").Value = 3130.4666
$ws.Range("K122").Value = 21518.715
$ws.Range("L122").Value = 9391.399800000001
$ws.Range("M122").Value = -19068.715
$ws.Range("N122").Value = -14291.3998

$ws.Range("H126").Value = 1900.5
$ws.Range("I126").Value = 1801
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 5403
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -2933
$ws.Range("N126").Value = -10940

$ws.Range("H136").Value = 3179.8064
$ws.Range("I136").Value = 3516.3044
$ws.Range("J136").Value = 2212.375
$ws.Range("K136").Value = 10548.9132
$ws.Range("L136").Value = 6637.125
$ws.Range("M136").Value = -7998.913199999999
$ws.Range("N136").Value = -11737.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 24000
$ws.Range("J97").Value = 24000
$ws.Range("L97").Value = 24000
$ws.Range("N97").Value = -25982

$ws.Range("H132").Value = 29517.65
$ws.Range("I132").Value = 40403.96
$ws.Range("J132").Value = 9300.214
$ws.Range("K132").Value = 121211.88
$ws.Range("L132").Value = 27900.642
$ws.Range("M132").Value = -118681.88
$ws.Range("N132").Value = -32960.642
